# Update column C ("Förändrad") from serial date 45182 (2023-09-13)
# to serial date 45184 (2023-09-15) for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -eq 45182) {
        $cell.Value2 = 45184
    }
}
